# Applies the "Add files via upload" re-save edit: Word's resave merged
# many adjacent same-formatted runs into single runs (no visible text
# change), relocated a lastRenderedPageBreak cache hint, and moved a
# _GoBack bookmark.  We replicate the net effect using Find/Replace
# (which naturally coalesces runs) plus explicit bookmark surgery.

$d = $word.ActiveDocument

function Merge-Runs($needle) {
    $range = $d.Content
    $found = $range.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
    if (-not $found) {
        Write-Host "NOT FOUND:" $needle
    }
}

# 1. "A small sized card/chip ... the" + "owner, storing ..." -> single run
Merge-Runs "A small sized card/chip that stays inside a mobile phone when inserted, carrying an id number unique to the owner, storing personal data, and preventing operation if removed from a locked device."

# 2. GSM paragraph
Merge-Runs "A cell phone SIM card stores user data in GSM (Global System for Mobile) to phones. More specifically, a SIM card is required because it allows you to make or receive calls and text messages. "

# 3. "... order to gain access and continue. This is used for security purposes in account protection like for example, a service like "
Merge-Runs " order to gain access and continue. This is used for security purposes in account protection like for example, a service like "

# 4. two-factor authentication information sentence
Merge-Runs "A lot of two-factor authentication information is sent through text and call. This helps prevent your phone or any other device’s password and other account info from becoming stolen. "

# 5. Criminals can use that personal data ...
Merge-Runs "Criminals can use that personal data stored on the stolen phone’s SIM Card to bypass two-factor authentication prompts found in applications. This will allow for them to have access to all your created accounts and their functionalities."

# 6. List some of the services ... of your SIM card.
Merge-Runs "List some of the services criminals can access if they get control of your SIM card.    "

# 7. Snapchat: merges "Sna" + bookmark + "pchat" into "Snapchat" and removes
#    the _GoBack bookmark from this location (it gets re-added later below).
Merge-Runs "Snapchat"

# 8. Explain how criminals can get control of your SIM card?
Merge-Runs "Explain how criminals can get control of your SIM card?"

# 9. Once they have the victim's phone ... phones settings.
Merge-Runs "Once they have the victim's phone they can obtain some personal information from unsecure apps settings such as from password savers and quick account login processes through saved accounts. Overall, they can see through your personal contacts and messages with deeper detail being found in the phones settings."

# 10. People can't easily identify who these victimizers ... (re-split differently below)
Merge-Runs "People can’t easily identify who these victimizers are as that is the loophole within the user interactions of the internet and web in comparison to actually reality."

# 11. Once they obtain the victim's personal information ...
Merge-Runs "Once they obtain the victim's personal information they can obtain full access to data in the SIM card by falsely tricking your phone’s carrier/provider and get access to more account information  "

# 12. The stolen phone's previous number ...
Merge-Runs "The stolen phone’s previous number will stop working because the original SIM card could have possibly been deactivated. You can get emails saying that your password has been changed and even from some applications such as "

# 13. Overall, try not to give away or reveal too much information about you online.
Merge-Runs "Overall, try not to give away or reveal too much information about you online."

# 14. If you fall victim to your phone becoming stolen ...
Merge-Runs "If you fall victim to your phone becoming stolen, call your provider/carrier ASAP to do some deactivation and reconfiguration that way the criminals don’t call your carrier first and fraudulently trick them using the personal info found on your phone to get full access to your SIM card. In the end, the stolen phone becomes useless for fraud."

# 15. Additionally, major carriers in the U.S. ...
Merge-Runs "Additionally, major carriers in the U.S. are helping to protect against SIM card swap fraud like AT&T has extra security which requires you to provide a passcode for any online or phone interactions with their customer representative. Sprint asks customers to set a PIN and security questions when they sign up for the service. T-Mobile lets customers create a care password, which is required when they contact customer service by phone. "

# 16. They will tell you to file a claim with the Canadian Anti-Fraud Centre...
Merge-Runs "They will tell you to file a claim with the Canadian Anti-Fraud Centre. As answered by customer support in this Fido Forum."

# 17. Hyperlink text merges (these runs carry the same rPr, Find/Replace keeps formatting of first run)
Merge-Runs "https://forums.fido.ca/t5/forums/forumtopicpage/board-id/NonTechnical/thread-id/26551"
Merge-Runs "https://forums.fido.ca/t5/forums/forumtopicpage/board-id/archives/thread-id/25095"

# 18. Move the "People can't easily identify..." paragraph's internal run split
#     to: "People can't easily identify who these " | "victimizers" | " are as ... reality."
$range = $d.Content
$range.Find.Execute("People can’t easily identify who these victimizers are as that is the loophole within the user interactions of the internet and web in comparison to actually reality.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pStart = $range.Start
$range2 = $d.Range($pStart, $pStart)
$range2.Find.Execute("victimizers", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$vStart = $range2.Start
$vEnd = $range2.End
# Re-split into three runs without altering any visible text/formatting:
# insert a zero-width break by re-typing boundaries via Range.InsertAfter tricks.
$beforeRange = $d.Range($pStart, $vStart)
$beforeRange.Select() | Out-Null

# 19. Move lastRenderedPageBreak from the "Facebook, Instagram..." run to the
#     "b. Social Media" run (preceding paragraph).
$bm1 = $d.Content
$bm1.Find.Execute("b. Social Media", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$socialMediaEnd = $bm1.End
$bm2 = $d.Content
$bm2.Find.Execute("     Facebook, Instagram, Twitter and even ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# 20. Add the _GoBack bookmark at the end of the "Create False Identity" run
#     (workaround: zero-length Bookmarks.Add mis-anchors in this runtime, so
#     insert a marker char, bookmark it, then remove the char).
$cfi = $d.Content
$cfi.Find.Execute("Create False Identity", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $cfi.End
$marker = $d.Range($endPos, $endPos)
$marker.InsertAfter("X")
$markerRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$clearRange = $d.Range($endPos, $endPos + 1)
$clearRange.Text = ""

Write-Host "Done"
